# Generate Report for Handoff
#
# The localization-status report re-sorts the three "file" rows (7, 8, 9)
# of every worksheet (Overview, zh-cn, de-de) by file name, and refreshes
# the status / handoff timestamps for the file that has just become
# "Ready for handoff" (fe721b5d...), which moves from row 7 to row 9.
#
#   old row7 (fe721b5d, In Translation)   -> new row9 (Ready for handoff, refreshed dates)
#   old row8 (5736df43, Ready for handoff) -> new row7 (unchanged)
#   old row9 (cbcdd771, Ready for handoff) -> new row8 (unchanged)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A7").Value = "5736df43-0ff9-48fb-9f0b-1f9845ced142.md"
$wsOverview.Range("B7").Value = "e2e\5736df43-0ff9-48fb-9f0b-1f9845ced142.md"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-08-21 00:48:45"

$wsOverview.Range("A8").Value = "cbcdd771-a78d-4d08-b66d-488e5b202a24.md"
$wsOverview.Range("B8").Value = "e2e\cbcdd771-a78d-4d08-b66d-488e5b202a24.md"
$wsOverview.Range("E8").Value = "Ready for handoff"
$wsOverview.Range("F8").Value = "Ready for handoff"
$wsOverview.Range("G8").Value = "2016-08-21 00:46:51"

$wsOverview.Range("A9").Value = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
$wsOverview.Range("B9").Value = "e2e\fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-08-21 00:51:30"

$overviewDisplay = @{
    '$B$7' = "e2e\5736df43-0ff9-48fb-9f0b-1f9845ced142.md"
    '$B$8' = "e2e\cbcdd771-a78d-4d08-b66d-488e5b202a24.md"
    '$B$9' = "e2e\fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
}
foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($overviewDisplay.ContainsKey($addr)) {
        $hl.TextToDisplay = $overviewDisplay[$addr]
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A7").Value = "5736df43-0ff9-48fb-9f0b-1f9845ced142.md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("G7").Value = "5736df43-0ff9-48fb-9f0b-1f9845ced142.3f16b4cde050102af2afd3e79794ce114dbd6be1.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-08-21 00:48:41"

$wsZhCn.Range("A8").Value = "cbcdd771-a78d-4d08-b66d-488e5b202a24.md"
$wsZhCn.Range("C8").Value = "Ready for handoff"
$wsZhCn.Range("G8").Value = "cbcdd771-a78d-4d08-b66d-488e5b202a24.b8e4142af020d03b283755bd354fcda2d644bedb.zh-cn.xlf"
$wsZhCn.Range("H8").Value = "2016-08-21 00:46:47"

$wsZhCn.Range("A9").Value = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("G9").Value = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.ee0f9686e6d1b7fa1b6ce94b379570f14e04ced2.zh-cn.xlf"
$wsZhCn.Range("H9").Value = "2016-08-21 00:51:26"

$zhcnDisplay = @{
    '$A$7' = "5736df43-0ff9-48fb-9f0b-1f9845ced142.md"
    '$A$8' = "cbcdd771-a78d-4d08-b66d-488e5b202a24.md"
    '$A$9' = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
}
foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($zhcnDisplay.ContainsKey($addr)) {
        $hl.TextToDisplay = $zhcnDisplay[$addr]
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A7").Value = "5736df43-0ff9-48fb-9f0b-1f9845ced142.md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("G7").Value = "5736df43-0ff9-48fb-9f0b-1f9845ced142.3f16b4cde050102af2afd3e79794ce114dbd6be1.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-08-21 00:48:45"

$wsDeDe.Range("A8").Value = "cbcdd771-a78d-4d08-b66d-488e5b202a24.md"
$wsDeDe.Range("C8").Value = "Ready for handoff"
$wsDeDe.Range("G8").Value = "cbcdd771-a78d-4d08-b66d-488e5b202a24.b8e4142af020d03b283755bd354fcda2d644bedb.de-de.xlf"
$wsDeDe.Range("H8").Value = "2016-08-21 00:46:51"

$wsDeDe.Range("A9").Value = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("G9").Value = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.ee0f9686e6d1b7fa1b6ce94b379570f14e04ced2.de-de.xlf"
$wsDeDe.Range("H9").Value = "2016-08-21 00:51:30"

$dedeDisplay = @{
    '$A$7' = "5736df43-0ff9-48fb-9f0b-1f9845ced142.md"
    '$A$8' = "cbcdd771-a78d-4d08-b66d-488e5b202a24.md"
    '$A$9' = "fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md"
}
foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($dedeDisplay.ContainsKey($addr)) {
        $hl.TextToDisplay = $dedeDisplay[$addr]
    }
}
